$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update UserName values (column A) for existing rows
$ws.Range("A2").Value = "sneakypythontestuser"
$ws.Range("A3").Value = "sneakypythontestuser"
$ws.Range("A4").Value = "SneakyPythonTestUser"
$ws.Range("A5").Value = "sn3akypythontestuser"

# Add two new rows of test data
$ws.Range("B6").Value = "000000aa"
$ws.Range("C6").Value = $false
$ws.Range("D6").Value = "The UserName is empty"

$ws.Range("A7").Value = "sneakypythontestuser"
$ws.Range("C7").Value = $false
$ws.Range("D7").Value = "The Password is empty"

# Column A width now needs to fit the longer usernames
$ws.Columns.Item(1).AutoFit() | Out-Null

# Update the selected cell / active cell shown when the file was saved
$ws.Range("D8").Select()
